# Data Update for 13-04-2020
# Fills in the "division" column (A) for the data rows that belong to each
# division group (Dhaka, Chattogram, Sylhet, Rangpur, Mymensingh, Barishal),
# mirroring the division name + formatting already present on the group's
# first/header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Dhaka division (rows 3-13) -------------------------------------------
# Header row 2 (A2) already carries "Dhaka" with style s=5; rows 3-13 only
# need the value -- the style is already correct (s=5), so no format copy
# is required here.
foreach ($r in 3..13) {
    $ws.Range("A$r").Value = "Dhaka"
}

# --- Chattogram division (rows 15-19), header row 14 -----------------------
$ws.Range("A14").Copy() | Out-Null
foreach ($r in 15..19) {
    $ws.Range("A$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("A$r").Value = "Chattogram"
}

# --- Sylhet division (rows 21-23), header row 20 ---------------------------
$ws.Range("A20").Copy() | Out-Null
foreach ($r in 21..23) {
    $ws.Range("A$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("A$r").Value = "Sylhet"
}

# --- Rangpur division (rows 25-28), header row 24 --------------------------
$ws.Range("A24").Copy() | Out-Null
foreach ($r in 25..28) {
    $ws.Range("A$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("A$r").Value = "Rangpur"
}

# --- Mymensingh division (rows 31-33), header row 30 -----------------------
$ws.Range("A30").Copy() | Out-Null
foreach ($r in 31..33) {
    $ws.Range("A$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("A$r").Value = "Mymensingh "
}

# --- Barishal division (rows 35-37), header row 34 -------------------------
$ws.Range("A34").Copy() | Out-Null
foreach ($r in 35..37) {
    $ws.Range("A$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("A$r").Value = "Barishal "
}

$excel.CutCopyMode = $false

# --- View state: scroll position + current selection -----------------------
$ws.Range("A124:XFD167").Select()
